$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 -----------------------------------------------------------------
# Date, submitter, rich-text model notes (partially bold) and time-spent.
# Copy the date formatting (built-in date number format) from an existing
# row so we reuse the workbook's existing style instead of minting a new one.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 43611

$ws.Range("B6").Value = "李希君"

$row6Text = "1. 表格merge: HASH_MAX = 200 WINDOW_SIZE = 5; `n2. categorical hash: factorization; `n3. feature selection simple`n4. 模型ensemble"
$ws.Range("C6").Value = $row6Text
$ws.Range("C6").WrapText = $true
$ws.Range("C6").Characters(86, 24).Font.Bold = $true
$tail6 = $ws.Range("C6").Characters(110, 14).Font
$tail6.Name = "Calibri"
$tail6.Size = 12

$ws.Range("D6").Value = 48

$ws.Rows.Item(6).RowHeight = 64

# --- Row 7 -----------------------------------------------------------------
# Date, submitter and plain-text model notes (no time-spent recorded).
$ws.Range("A5").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 43611

$ws.Range("B7").Value = "李希君"

$row7Text = "1. 表格merge: HASH_MAX = 200 WINDOW_SIZE = 5; `n2. categorical hash: factorization`n3. 模型ensemble"
$ws.Range("C7").Value = $row7Text
$ws.Range("C7").WrapText = $true

$ws.Rows.Item(7).RowHeight = 48

# --- Selection mirrors the author's last-edited cell -----------------------
$null = $ws.Range("D7").Select()
